$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 193.33333
$ws.Range("I2").Value = 164.28572
$ws.Range("J2").Value = 295
$ws.Range("K2").Value = 164.28572
$ws.Range("L2").Value = 295
$ws.Range("M2").Value = -51.28572
$ws.Range("N2").Value = -521
$ws.Range("H17").Value = 1111.1111
$ws.Range("I17").Value = 500
$ws.Range("J17").Value = 1147.0588
$ws.Range("K17").Value = 1500
$ws.Range("L17").Value = 3441.1764
$ws.Range("M17").Value = -1332
$ws.Range("N17").Value = -3777.1764
$ws.Range("H125").Value = 714.0323
$ws.Range("I125").Value = 899.8333
$ws.Range("J125").Value = 669.4400000000001
$ws.Range("K125").Value = 8098.4997
$ws.Range("L125").Value = 6024.960000000001
$ws.Range("M125").Value = -5638.4997
$ws.Range("N125").Value = -10944.96
$ws.Range("H129").Value = 713.76624
$ws.Range("I129").Value = 367.33334
$ws.Range("J129").Value = 743.04224
$ws.Range("K129").Value = 1102.00002
$ws.Range("L129").Value = 2229.12672
$ws.Range("M129").Value = 3897.99998
$ws.Range("N129").Value = -12229.12672
$ws.Range("H130").Value = 23065
$ws.Range("J130").Value = 23065
$ws.Range("L130").Value = 23065
$ws.Range("N130").Value = -33105
$ws.Range("H131").Value = 1833.775
$ws.Range("I131").Value = 1251
$ws.Range("K131").Value = 3753
$ws.Range("M131").Value = 1287
$ws.Range("H132").Value = 1777.4642
$ws.Range("I132").Value = 1627.9636
$ws.Range("J132").Value = 10000
$ws.Range("K132").Value = 4883.8908
$ws.Range("L132").Value = 30000
$ws.Range("M132").Value = -2353.8908
$ws.Range("N132").Value = -35060
$ws.Range("H133").Value = 39568.57
$ws.Range("J133").Value = 39568.57
$ws.Range("L133").Value = 39568.57
$ws.Range("N133").Value = -49688.57
$ws.Range("H134").Value = 38000
$ws.Range("J134").Value = 38000
$ws.Range("L134").Value = 38000
$ws.Range("N134").Value = -48140
$ws.Range("H138").Value = 1797.5784
$ws.Range("I138").Value = 967.0741
$ws.Range("J138").Value = 2198
$ws.Range("K138").Value = 2901.2223
$ws.Range("L138").Value = 6594
$ws.Range("M138").Value = 2238.7777
$ws.Range("N138").Value = -16874
$ws.Range("H140").Value = 33880
$ws.Range("J140").Value = 33880
$ws.Range("L140").Value = 33880
$ws.Range("N140").Value = -44240

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1437.5
$ws.Range("I122").Value = 1200
$ws.Range("J122").Value = 1833.3334
$ws.Range("K122").Value = 3600
$ws.Range("L122").Value = 5500.0002
$ws.Range("M122").Value = -1150
$ws.Range("N122").Value = -10400.0002

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 473.3158
$ws.Range("I94").Value = 415.2143
$ws.Range("J94").Value = 636
$ws.Range("K94").Value = 415.2143
$ws.Range("L94").Value = 636
$ws.Range("M94").Value = 35.78570000000002
$ws.Range("N94").Value = -1538

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2959.9773
$ws.Range("I31").Value = 2050.8928
$ws.Range("J31").Value = 4550.875
$ws.Range("K31").Value = 2050.8928
$ws.Range("L31").Value = 4550.875
$ws.Range("M31").Value = -1755.8928
$ws.Range("N31").Value = -5140.875
$ws.Range("H34").Value = 2959.9773
$ws.Range("I34").Value = 2050.8928
$ws.Range("J34").Value = 4550.875
$ws.Range("K34").Value = 2050.8928
$ws.Range("L34").Value = 4550.875
$ws.Range("M34").Value = -1848.8928
$ws.Range("N34").Value = -4954.875
$ws.Range("H110").Value = 57415.332
$ws.Range("J110").Value = 57415.332
$ws.Range("L110").Value = 57415.332
$ws.Range("N110").Value = -65595.33199999999
$ws.Range("H111").Value = 28106.25
$ws.Range("J111").Value = 28106.25
$ws.Range("L111").Value = 28106.25
$ws.Range("N111").Value = -36286.25
$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").ClearContents()
$ws.Range("H116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("N116").ClearContents()
$ws.Range("H119").Value = 48925
$ws.Range("J119").Value = 48925
$ws.Range("L119").Value = 48925
$ws.Range("N119").Value = -58601

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 5400
$ws.Range("J9").Value = 5400
$ws.Range("L9").Value = 16200
$ws.Range("N9").Value = -16648
$ws.Range("H41").Value = 251.42857
$ws.Range("I41").Value = 152
$ws.Range("J41").Value = 500
$ws.Range("K41").Value = 456
$ws.Range("L41").Value = 1500
$ws.Range("M41").Value = -118
$ws.Range("N41").Value = -2176

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2071.3333
$ws.Range("I102").Value = 1956
$ws.Range("J102").Value = 2475
$ws.Range("K102").Value = 1956
$ws.Range("L102").Value = 2475
$ws.Range("M102").Value = -334
$ws.Range("N102").Value = -5719
$ws.Range("H126").Value = 4764176
$ws.Range("I126").Value = 7144475
$ws.Range("J126").Value = 3577.4285
$ws.Range("K126").Value = 21433425
$ws.Range("L126").Value = 10732.2855
$ws.Range("M126").Value = -21430955
$ws.Range("N126").Value = -15672.2855
$ws.Range("H132").Value = 4693.7036
$ws.Range("I132").Value = 5168.5
$ws.Range("J132").Value = 4003.0908
$ws.Range("K132").Value = 15505.5
$ws.Range("L132").Value = 12009.2724
$ws.Range("M132").Value = -12975.5
$ws.Range("N132").Value = -17069.2724

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 271.8125
$ws.Range("I55").Value = 149.33333
$ws.Range("J55").Value = 345.3
$ws.Range("K55").Value = 149.33333
$ws.Range("L55").Value = 345.3
$ws.Range("M55").Value = 23.66667000000001
$ws.Range("N55").Value = -691.3
$ws.Range("H93").Value = 5121.3335
$ws.Range("I93").Value = 7976.0713
$ws.Range("J93").Value = 1124.7
$ws.Range("K93").Value = 7976.0713
$ws.Range("L93").Value = 1124.7
$ws.Range("M93").Value = -6728.0713
$ws.Range("N93").Value = -3620.7
$ws.Range("H132").Value = 2302.83
$ws.Range("I132").Value = 1693.884
$ws.Range("J132").Value = 3658.2258
$ws.Range("K132").Value = 5081.652
$ws.Range("L132").Value = 10974.6774
$ws.Range("M132").Value = -2551.652
$ws.Range("N132").Value = -16034.6774

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1991.4131
$ws.Range("I132").Value = 956.6070999999999
$ws.Range("J132").Value = 3601.111
$ws.Range("K132").Value = 2869.8213
$ws.Range("L132").Value = 10803.333
$ws.Range("M132").Value = -339.8212999999996
$ws.Range("N132").Value = -15863.333
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()

Write-Output "Applied all changes"